# Append three new daily rows (235-237) to the Quantities sheet, mirroring
# the last existing row (234) with column A (date) incremented by one day
# per new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 234
$newRowCount = 3
$lastCol = 10   # column J

for ($i = 1; $i -le $newRowCount; $i++) {
    $r = $lastRow + $i

    # Copy the number format/style of column A from the last existing row so
    # the new date cell reuses the same cellXf (style index) instead of
    # Excel creating a brand-new duplicate style when the value is set.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    # Column A: date serial number, incrementing by one day per new row.
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($lastRow, 1).Value2 + $i

    # Columns B..J: identical values to the last existing row.
    for ($c = 2; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($lastRow, $c).Value2
    }
}

$excel.CutCopyMode = 0
